$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "9/9/2018, 7:34:16 AM"
$ws.Range("B2").Value = "svmasdhruthi@gmail.com"
$ws.Range("C2").Value = "Chatbots in educational system"
$ws.Range("D2").Value = "Ms. MADHURA J."
$ws.Range("E2").Value = "Chatbots in educational system"
$ws.Range("F2").Value = "Ms. MADHURA J."
$ws.Range("G2").Value = "Detecting fraud apps using sentiment analysis"
$ws.Range("H2").Value = "Mr. SURESHKUMAR M."
$ws.Range("I2").Value = "Web based library management system using angular and springboot"
$ws.Range("J2").Value = "Ms. CHANDRAKALA B.M."
$ws.Range("K2").Value = "1DS15IS026"
$ws.Range("L2").Value = "1DS15IS046"
$ws.Range("M2").Value = "1DS15IS053"
$ws.Range("N2").Value = "1DS15IS057"
